$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 265.83334
$ws.Range("I33").Value = 280.9091
$ws.Range("K33").Value = 280.9091
$ws.Range("M33").Value = -51.90910000000002
$ws.Range("H43").Value = 1744.4615
$ws.Range("I43").Value = 1228
$ws.Range("J43").Value = 2067.25
$ws.Range("K43").Value = 1228
$ws.Range("L43").Value = 2067.25
$ws.Range("M43").Value = -1159
$ws.Range("N43").Value = -2205.25
$ws.Range("H74").Value = 5267890.5
$ws.Range("I74").Value = 10003950
$ws.Range("K74").Value = 10003950
$ws.Range("M74").Value = -10003014
$ws.Range("H76").Value = 4369.9
$ws.Range("I76").Value = 3159.8
$ws.Range("J76").Value = 5580
$ws.Range("K76").Value = 3159.8
$ws.Range("L76").Value = 5580
$ws.Range("M76").Value = -2844.8
$ws.Range("N76").Value = -6210
$ws.Range("H77").Value = 5267890.5
$ws.Range("I77").Value = 10003950
$ws.Range("K77").Value = 50019750
$ws.Range("M77").Value = -50015070
$ws.Range("H79").Value = 4369.9
$ws.Range("I79").Value = 3159.8
$ws.Range("J79").Value = 5580
$ws.Range("K79").Value = 3159.8
$ws.Range("L79").Value = 5580
$ws.Range("M79").Value = -2067.8
$ws.Range("N79").Value = -7764
$ws.Range("H100").Value = 33334784
$ws.Range("I100").Value = 40001520
$ws.Range("J100").Value = 1106
$ws.Range("K100").Value = 40001520
$ws.Range("L100").Value = 1106
$ws.Range("M100").Value = -40000979
$ws.Range("N100").Value = -2188
$ws.Range("H113").Value = 6839.5
$ws.Range("I113").Value = 3118.3333
$ws.Range("K113").Value = 3118.3333
$ws.Range("M113").Value = 135.6667000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4834.8887
$ws.Range("I32").Value = 3785.68
$ws.Range("J32").Value = 17950
$ws.Range("K32").Value = 3785.68
$ws.Range("L32").Value = 17950
$ws.Range("M32").Value = -3498.68
$ws.Range("N32").Value = -18524
$ws.Range("H63").Value = 7698649.5
$ws.Range("I63").Value = 19789810
$ws.Range("J63").Value = 4274.4546
$ws.Range("K63").Value = 19789810
$ws.Range("L63").Value = 4274.4546
$ws.Range("M63").Value = -19789124
$ws.Range("N63").Value = -5646.4546
$ws.Range("H66").Value = 7698649.5
$ws.Range("I66").Value = 19789810
$ws.Range("J66").Value = 4274.4546
$ws.Range("K66").Value = 98949050
$ws.Range("L66").Value = 21372.273
$ws.Range("M66").Value = -98945618
$ws.Range("N66").Value = -28236.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 46600
$ws.Range("J82").Value = 46600
$ws.Range("L82").Value = 46600
$ws.Range("N82").Value = -47366
$ws.Range("H85").Value = 46600
$ws.Range("J85").Value = 46600
$ws.Range("L85").Value = 46600
$ws.Range("N85").Value = -49252
$ws.Range("H130").Value = 41862.5
$ws.Range("J130").Value = 41862.5
$ws.Range("L130").Value = 41862.5
$ws.Range("N130").Value = -51902.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13337659
$ws.Range("I99").Value = 20001688
$ws.Range("J99").Value = 9600
$ws.Range("K99").Value = 20001688
$ws.Range("L99").Value = 9600
$ws.Range("M99").Value = -20000190
$ws.Range("N99").Value = -12596
$ws.Range("H122").Value = 5106.25
$ws.Range("I122").Value = 1555.5
$ws.Range("J122").Value = 8657
$ws.Range("K122").Value = 4666.5
$ws.Range("L122").Value = 25971
$ws.Range("M122").Value = -2216.5
$ws.Range("N122").Value = -30871
$ws.Range("H126").Value = 13337659
$ws.Range("I126").Value = 20001688
$ws.Range("J126").Value = 9600
$ws.Range("K126").Value = 60005064
$ws.Range("L126").Value = 28800
$ws.Range("M126").Value = -60002594
$ws.Range("N126").Value = -33740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 16094.462
$ws.Range("J80").Value = 17185.666
$ws.Range("L80").Value = 51556.99800000001
$ws.Range("N80").Value = -53428.99800000001
$ws.Range("H83").Value = 16094.462
$ws.Range("J83").Value = 17185.666
$ws.Range("L83").Value = 154670.994
$ws.Range("N83").Value = -164030.994
$ws.Range("H122").Value = 3162.4
$ws.Range("I122").Value = 1112.2858
$ws.Range("J122").Value = 3674.9285
$ws.Range("K122").Value = 10010.5722
$ws.Range("L122").Value = 33074.3565
$ws.Range("M122").Value = -7560.572200000001
$ws.Range("N122").Value = -37974.3565
$ws.Range("H123").Value = 3951.7273
$ws.Range("I123").Value = 3846.9
$ws.Range("K123").Value = 11540.7
$ws.Range("M123").Value = -9090.700000000001
$ws.Range("H124").Value = 4306.357
$ws.Range("I124").Value = 2581.6667
$ws.Range("J124").Value = 5599.875
$ws.Range("K124").Value = 7745.000100000001
$ws.Range("L124").Value = 16799.625
$ws.Range("M124").Value = -2835.000100000001
$ws.Range("N124").Value = -26619.625
$ws.Range("H125").Value = 5406
$ws.Range("I125").Value = 2015
$ws.Range("J125").Value = 7666.6665
$ws.Range("K125").Value = 6045
$ws.Range("L125").Value = 22999.9995
$ws.Range("M125").Value = -1125
$ws.Range("N125").Value = -32839.99950000001
$ws.Range("H131").Value = 709.8200000000001
$ws.Range("I131").Value = 265.10526
$ws.Range("J131").Value = 814.1358
$ws.Range("K131").Value = 795.3157799999999
$ws.Range("L131").Value = 2442.4074
$ws.Range("M131").Value = 4244.68422
$ws.Range("N131").Value = -12522.4074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6472.1665
$ws.Range("I70").Value = 5812.222
$ws.Range("K70").Value = 5812.222
$ws.Range("M70").Value = -5542.222
$ws.Range("H73").Value = 6472.1665
$ws.Range("I73").Value = 5812.222
$ws.Range("K73").Value = 5812.222
$ws.Range("M73").Value = -4876.222
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 70780
$ws.Range("J130").Value = 70780
$ws.Range("L130").Value = 70780
$ws.Range("N130").Value = -80820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3445.4348
$ws.Range("J7").Value = 5685.7144
$ws.Range("L7").Value = 5685.7144
$ws.Range("N7").Value = -5909.7144
$ws.Range("H16").Value = 1669.125
$ws.Range("I16").Value = 1478.7142
$ws.Range("J16").Value = 3002
$ws.Range("K16").Value = 1478.7142
$ws.Range("L16").Value = 3002
$ws.Range("M16").Value = -1308.7142
$ws.Range("N16").Value = -3342
$ws.Range("H40").Value = 7843.0713
$ws.Range("I40").Value = 7129
$ws.Range("K40").Value = 7129
$ws.Range("M40").Value = -6993
$ws.Range("H126").Value = 3445.4348
$ws.Range("J126").Value = 5685.7144
$ws.Range("L126").Value = 17057.1432
$ws.Range("N126").Value = -21997.1432
$ws.Range("H133").Value = 28890
$ws.Range("J133").Value = 28890
$ws.Range("L133").Value = 28890
$ws.Range("N133").Value = -33950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 22315.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 22315.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 22315.5
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -23355.5
$ws.Range("H81").Value = 4000
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 4000
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608
$ws.Range("H126").Value = 508599.25
$ws.Range("I126").Value = 1311.5
$ws.Range("J126").Value = 1184982.9
$ws.Range("K126").Value = 3934.5
$ws.Range("L126").Value = 3554948.7
$ws.Range("M126").Value = -1464.5
$ws.Range("N126").Value = -3559888.7
$ws.Range("H132").Value = 10418883
$ws.Range("I132").Value = 1605.1904
$ws.Range("J132").Value = 30306412
$ws.Range("K132").Value = 4815.5712
$ws.Range("L132").Value = 90919236
$ws.Range("M132").Value = -2285.5712
$ws.Range("N132").Value = -90924296
